# Fruta / hortaliza, semanal
# Insert two new daily-price rows (date serial 44636 = 2022-03-16)
# above the existing row 20, pushing the previous rows 20-25 down to 22-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 20 (shifts old rows 20-25 down to 22-27).
$ws.Range("A20:T21").Insert()

# New row 20: Tuna, Especial, 2022-03-16 (serial 44636)
$ws.Range("A20").Value = 8
$ws.Range("B20").Value = "Terminal La Palmera de La Serena"
$ws.Range("C20").Value = "Coquimbo"
$ws.Range("D20").Value = 44636
$ws.Range("E20").Value = 4
$ws.Range("F20").Value = "Fruta"
$ws.Range("G20").Value = 100107
$ws.Range("H20").Value = "Otros"
$ws.Range("I20").Value = 100107011
$ws.Range("J20").Value = "Tuna"
$ws.Range("K20").Value = "Sin especificar"
$ws.Range("L20").Value = "Especial"
$ws.Range("M20").Value = 240
$ws.Range("N20").Value = 14000
$ws.Range("O20").Value = 15000
$ws.Range("P20").Value = 14500
$ws.Range("Q20").Value = "$/caja 18 kilos"
$ws.Range("R20").Value = "Provincia de Limarí"
$ws.Range("S20").Value = 806
$ws.Range("T20").Value = 18

# New row 21: Tuna, Primera, 2022-03-16 (serial 44636)
$ws.Range("A21").Value = 8
$ws.Range("B21").Value = "Terminal La Palmera de La Serena"
$ws.Range("C21").Value = "Coquimbo"
$ws.Range("D21").Value = 44636
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = "Fruta"
$ws.Range("G21").Value = 100107
$ws.Range("H21").Value = "Otros"
$ws.Range("I21").Value = 100107011
$ws.Range("J21").Value = "Tuna"
$ws.Range("K21").Value = "Sin especificar"
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 200
$ws.Range("N21").Value = 10000
$ws.Range("O21").Value = 11000
$ws.Range("P21").Value = 10500
$ws.Range("Q21").Value = "$/caja 18 kilos"
$ws.Range("R21").Value = "Provincia de Limarí"
$ws.Range("S21").Value = 583
$ws.Range("T21").Value = 18

# The newly inserted rows copy formatting from the row above/below; make sure
# the date columns keep the date number format used throughout column D.
$ws.Range("D20:D21").NumberFormat = $ws.Range("D22").NumberFormat
